$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Range G1 Width (px):" $ws.Range("G1").Width
Write-Host "Range F1 Width (px):" $ws.Range("F1").Width
Write-Host "Range B1 Width (px):" $ws.Range("B1").Width
